$wb = $excel.ActiveWorkbook

# --- DatosCuenta (sheet1): rename Smoke QA record to Smoke PreProd record, bump numeric values ---
$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsCuenta.Range("A2").Value = "SmokPreProdOcho"
$wsCuenta.Range("B2").Value = "SmokeNamePreProdOcho"
$wsCuenta.Range("C2").Value = 27100116
$wsCuenta.Range("D2").Value = 118

# --- DatosHogar (sheet2): bump numeric value ---
$wsHogar = $wb.Worksheets.Item("DatosHogar")
$wsHogar.Range("A2").Value = 637

# --- DatosMotor (sheet3): rename plate/engine/chassis codes from SMA018 to SMP019 ---
$wsMotor = $wb.Worksheets.Item("DatosMotor")
$wsMotor.Range("A2").Value = "SMP019"
$wsMotor.Range("B2").Value = "ABC12SSMP019"
$wsMotor.Range("C2").Value = "ZAZ123SSMP019"

# --- DatosAP (sheet4): bump numeric value ---
$wsAP = $wb.Worksheets.Item("DatosAP")
$wsAP.Range("A2").Value = 21200119

# --- Update the selected cell on each sheet to match the new active cells ---
$wsCuenta.Range("D2").Select()
$wsHogar.Range("A2").Select()
$wsMotor.Range("C3").Select()
$wsAP.Range("A3").Select()
